$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at position 8 ---------------------------------
# This pushes the old row 8 (the "total" row, P8/Q8 = 50.5) down to row 9
# and the old row 9 (footer row) down to row 10. Excel automatically
# shifts the merged-cell references and keeps the existing style indices
# for the rows that moved.
$ws.Rows.Item(8).Insert()

# --- 2. Clone row 7's cell formatting onto the new row 8 ----------------
# Row 7 is the existing "item" row (DOXIRAZOL ...). The new row 8 becomes
# a second item row (FLECTOR ...), so copy formats only (not values).
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)

# --- 3. Fill in the new item row's values --------------------------------
$ws.Cells.Item(8, 1).Value = 2                         # A8: item number
$ws.Cells.Item(8, 3).Value = "FLECTOR 50MG 30 CAPS"    # C8: item name
$ws.Cells.Item(8, 8).Value = "0:0"                     # H8
$ws.Cells.Item(8, 12).Value = "'1"                     # L8 (force text "1")
$ws.Cells.Item(8, 14).Value = "87.00"                  # N8
$ws.Cells.Item(8, 16).Value = "'87.0000"               # P8 (force text)
$ws.Cells.Item(8, 17).Value = "1:0"                    # Q8

# Re-paste the original number formats on the two cells where we had to
# force a leading apostrophe (quote-prefix resets the style), so they end
# up with the same number format as their row-7 counterparts again.
$ws.Range("L7").Copy()
$ws.Range("L8").PasteSpecial(-4122)
$ws.Range("P7").Copy()
$ws.Range("P8").PasteSpecial(-4122)

# --- 4. Re-create the merged cells for the new row 8 ---------------------
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

# --- 5. Row heights --------------------------------------------------------
$ws.Rows.Item(8).RowHeight = 24.75   # new item row
$ws.Rows.Item(9).RowHeight = 25.5    # total row (was 24.75 before insert)

# --- 6. Update the total in row 9 (was row 8) -----------------------------
$ws.Cells.Item(9, 16).Value = 137.5  # P9: 50.5 + 87.00

# --- 7. Update the footer timestamp in row 10 (was row 9) ----------------
$ws.Cells.Item(10, 1).Value = "Sunday, 28 September, 2025 9:39 AM"
